$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.15514874458313
$ws.Range("B1").Value = 1.728742957115173
$ws.Range("C1").Value = 7.039756298065186
$ws.Range("D1").Value = 2.65102744102478
$ws.Range("E1").Value = 1.46392285823822
